$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1525.2727
$ws.Range("I40").Value = 1199.8334
$ws.Range("J40").Value = 1647.3125
$ws.Range("K40").Value = 1199.8334
$ws.Range("L40").Value = 1647.3125
$ws.Range("M40").Value = -1024.8334
$ws.Range("N40").Value = -1997.3125

$ws.Range("H75").Value = 34993
$ws.Range("J75").Value = 34993
$ws.Range("L75").Value = 34993
$ws.Range("N75").Value = -36865

$ws.Range("H78").Value = 34993
$ws.Range("J78").Value = 34993
$ws.Range("L78").Value = 104979
$ws.Range("N78").Value = -114339

$ws.Range("H98").Value = 1318.6666
$ws.Range("I98").Value = 1371.6666
$ws.Range("J98").Value = 1106.6666
$ws.Range("K98").Value = 1371.6666
$ws.Range("L98").Value = 1106.6666
$ws.Range("M98").Value = 126.3334
$ws.Range("N98").Value = -4102.6666

$ws.Range("H122").Value = 1318.6666
$ws.Range("I122").Value = 1371.6666
$ws.Range("J122").Value = 1106.6666
$ws.Range("K122").Value = 4114.9998
$ws.Range("L122").Value = 3319.9998
$ws.Range("M122").Value = -1664.9998
$ws.Range("N122").Value = -8219.9998

$ws.Range("H125").Value = 3454
$ws.Range("I125").Value = 6877.3335
$ws.Range("J125").Value = 1400
$ws.Range("K125").Value = 61896.0015
$ws.Range("L125").Value = 12600
$ws.Range("M125").Value = -59436.0015
$ws.Range("N125").Value = -17520

$ws.Range("H137").Value = 66667336
$ws.Range("J137").Value = 200000000
$ws.Range("L137").Value = 600000000
$ws.Range("N137").Value = -600005100

$ws.Range("H141").Value = 1514.4231
$ws.Range("I141").Value = 1208.3334
$ws.Range("J141").Value = 2800
$ws.Range("K141").Value = 3625.0002
$ws.Range("L141").Value = 8400
$ws.Range("M141").Value = 1554.9998
$ws.Range("N141").Value = -18760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6874.5527
$ws.Range("I32").Value = 7124.776
$ws.Range("J32").Value = 6068.278
$ws.Range("K32").Value = 7124.776
$ws.Range("L32").Value = 6068.278
$ws.Range("M32").Value = -6837.776
$ws.Range("N32").Value = -6642.278

$ws.Range("H74").Value = 62508104
$ws.Range("I74").Value = 125005704
$ws.Range("J74").Value = 10507
$ws.Range("K74").Value = 125005704
$ws.Range("L74").Value = 10507
$ws.Range("M74").Value = -125004830
$ws.Range("N74").Value = -12255

$ws.Range("H77").Value = 62508104
$ws.Range("I77").Value = 125005704
$ws.Range("J77").Value = 10507
$ws.Range("K77").Value = 625028520
$ws.Range("L77").Value = 52535
$ws.Range("M77").Value = -625024152
$ws.Range("N77").Value = -61271

$ws.Range("H107").Value = 29000
$ws.Range("J107").Value = 29000
$ws.Range("L107").Value = 29000
$ws.Range("N107").Value = -36680

$ws.Range("H132").Value = 9262121
$ws.Range("I132").Value = 13891047
$ws.Range("J132").Value = 4269.5557
$ws.Range("K132").Value = 41673141
$ws.Range("L132").Value = 12808.6671
$ws.Range("M132").Value = -41670611
$ws.Range("N132").Value = -17868.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 42097
$ws.Range("J69").Value = 42097
$ws.Range("L69").Value = 42097
$ws.Range("N69").Value = -43719

$ws.Range("H72").Value = 42097
$ws.Range("J72").Value = 42097
$ws.Range("L72").Value = 126291
$ws.Range("N72").Value = -134403

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7579738.5
$ws.Range("I31").Value = 4137.6587
$ws.Range("J31").Value = 111112950
$ws.Range("K31").Value = 4137.6587
$ws.Range("L31").Value = 111112950
$ws.Range("M31").Value = -3842.6587
$ws.Range("N31").Value = -111113540

$ws.Range("H34").Value = 7579738.5
$ws.Range("I34").Value = 4137.6587
$ws.Range("J34").Value = 111112950
$ws.Range("K34").Value = 4137.6587
$ws.Range("L34").Value = 111112950
$ws.Range("M34").Value = -3935.6587
$ws.Range("N34").Value = -111113354

$ws.Range("H68").Value = 25863.334
$ws.Range("J68").Value = 28636
$ws.Range("L68").Value = 28636
$ws.Range("N68").Value = -30134

$ws.Range("H71").Value = 25863.334
$ws.Range("J71").Value = 28636
$ws.Range("L71").Value = 85908
$ws.Range("N71").Value = -93396

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 16666812
$ws.Range("I7").Value = 25000120
$ws.Range("J7").Value = 197.5
$ws.Range("K7").Value = 75000360
$ws.Range("L7").Value = 592.5
$ws.Range("M7").Value = -75000248
$ws.Range("N7").Value = -816.5

$ws.Range("H39").Value = 510.18182
$ws.Range("J39").Value = 510.18182
$ws.Range("L39").Value = 1530.54546
$ws.Range("N39").Value = -2118.54546

$ws.Range("H82").Value = 7836.1763
$ws.Range("I82").Value = 706.5
$ws.Range("J82").Value = 8786.799999999999
$ws.Range("K82").Value = 2119.5
$ws.Range("L82").Value = 26360.4
$ws.Range("M82").Value = -1713.5
$ws.Range("N82").Value = -27172.4

$ws.Range("H85").Value = 7836.1763
$ws.Range("I85").Value = 706.5
$ws.Range("J85").Value = 8786.799999999999
$ws.Range("K85").Value = 2119.5
$ws.Range("L85").Value = 26360.4
$ws.Range("M85").Value = -715.5
$ws.Range("N85").Value = -29168.4

$ws.Range("H113").Value = 1018.9091
$ws.Range("I113").Value = 549.3333
$ws.Range("J113").Value = 1344
$ws.Range("K113").Value = 1647.9999
$ws.Range("L113").Value = 4032
$ws.Range("M113").Value = 522.0001
$ws.Range("N113").Value = -8372

$ws.Range("H121").Value = 792.8929000000001
$ws.Range("I121").Value = 183.75
$ws.Range("J121").Value = 1036.55
$ws.Range("K121").Value = 551.25
$ws.Range("L121").Value = 3109.65
$ws.Range("M121").Value = 758.75
$ws.Range("N121").Value = -5729.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 30654.334
$ws.Range("I70").Value = 50000
$ws.Range("J70").Value = 20981.5
$ws.Range("K70").Value = 50000
$ws.Range("L70").Value = 20981.5
$ws.Range("M70").Value = -49730
$ws.Range("N70").Value = -21521.5

$ws.Range("H73").Value = 30654.334
$ws.Range("I73").Value = 50000
$ws.Range("J73").Value = 20981.5
$ws.Range("K73").Value = 50000
$ws.Range("L73").Value = 20981.5
$ws.Range("M73").Value = -49064
$ws.Range("N73").Value = -22853.5

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12453.9
$ws.Range("I62").Value = 5933.3335
$ws.Range("J62").Value = 15248.429
$ws.Range("K62").Value = 5933.3335
$ws.Range("L62").Value = 15248.429
$ws.Range("M62").Value = -5309.3335
$ws.Range("N62").Value = -16496.429

$ws.Range("H65").Value = 12453.9
$ws.Range("I65").Value = 5933.3335
$ws.Range("J65").Value = 15248.429
$ws.Range("K65").Value = 29666.6675
$ws.Range("L65").Value = 76242.145
$ws.Range("M65").Value = -26546.6675
$ws.Range("N65").Value = -82482.145

$ws.Range("H122").Value = 2682.7273
$ws.Range("I122").Value = 2682.7273
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8048.1819
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5598.1819
$ws.Range("N122").ClearContents()
